# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from 2023-10-05 (45204) to 2023-10-08 (45207).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
